$d = $word.ActiveDocument

# --- Update the date cells (column 1) to include the year ---
$t0 = $d.Tables.Item(1)

$r2 = $t0.Cell(2, 1).Range
$r2.MoveEnd(1, -1)
$r2.Text = "December 12th 2021 "

$r3 = $t0.Cell(3, 1).Range
$r3.MoveEnd(1, -1)
$r3.Text = "December 27th 2021"

$r4 = $t0.Cell(4, 1).Range
$r4.MoveEnd(1, -1)
$r4.Text = "December 29th 2021"

$r5 = $t0.Cell(5, 1).Range
$r5.MoveEnd(1, -1)
$r5.Text = "December 30th 2021"

$r6 = $t0.Cell(6, 1).Range
$r6.MoveEnd(1, -1)
$r6.Text = "December 31st 2021"

# --- Append a new paragraph to the last (Zubair) cell of the Dec 31st row ---
$t = $d.Tables.Item(1)
$lastRowIdx = $t.Rows.Count
$zubairCell = $t.Cell($lastRowIdx, 3)
$r = $zubairCell.Range
$r.Collapse(0)
$cr = [string][char]13
$r.InsertAfter($cr + "- Continued working on flow chart")

# --- Add a new row for January 1st 2022 ---
$t2 = $d.Tables.Item(1)
$t2.Rows.Add() | Out-Null
$newRowIdx = $t2.Rows.Count
$newRow = $t2.Rows.Item($newRowIdx)
$newRow.HeadingFormat = 0
$t2.Cell($newRowIdx, 1).Range.Text = "January 1st 2022"
$t2.Cell($newRowIdx, 3).Range.Text = "- Continued working on flow chart"

Write-Output "All edits applied"
